$wb = $excel.ActiveWorkbook

# --- Generate Report for handback ---
# For each localized-language sheet (zh-cn, de-de), fill in:
#   - Status (col B)                -> "Handed back: in sync with en-us"
#   - Latest Target File (col E)    -> same file name/link as Source File Name (col A)
#   - Latest Handback File (col F)  -> same file name/link as Latest Handoff File (col C)
#   - Latest Handback DateTime (col G) -> timestamp of the handback

$langSheets = @(
    @{ Name = "zh-cn"; HandbackTime = "2016-01-08 19:14:08"; MdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/ff527b8315e97b529250d8388bde38c31c232910/e2e/3d5df28b-542e-4401-841a-55279b6c8572.md"; XlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/90fd955952248cbbfec0f0a319da67d886247192/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/3d5df28b-542e-4401-841a-55279b6c8572.60c8402e9e0df38573a4c8226aa6d78291a5572a.zh-cn.xlf"; XlfName = "3d5df28b-542e-4401-841a-55279b6c8572.60c8402e9e0df38573a4c8226aa6d78291a5572a.zh-cn.xlf" },
    @{ Name = "de-de"; HandbackTime = "2016-01-08 19:14:24"; MdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/ff527b8315e97b529250d8388bde38c31c232910/e2e/3d5df28b-542e-4401-841a-55279b6c8572.md"; XlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/20db7265be7b2c2d375b593b4e94b5ca379b6ebe/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/3d5df28b-542e-4401-841a-55279b6c8572.60c8402e9e0df38573a4c8226aa6d78291a5572a.de-de.xlf"; XlfName = "3d5df28b-542e-4401-841a-55279b6c8572.60c8402e9e0df38573a4c8226aa6d78291a5572a.de-de.xlf" }
)

$mdFileName = "3d5df28b-542e-4401-841a-55279b6c8572.md"

# The "Overview" sheet rolls up each language's Status in columns B/C, so it
# picks up the same "Handed back" wording shown on the per-language sheets.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-us"
$overview.Range("C2").Value = "Handed back: in sync with en-us"

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # Status -> handed back
    $ws.Range("B2").Value = "Handed back: in sync with en-us"

    # Latest Target File (E2) - mirrors Source File Name (A2) hyperlink
    $ws.Hyperlinks.Add($ws.Range("E2"), $lang.MdAddress, "", "", $mdFileName)
    $ws.Range("E2").Font.Underline = $true
    $ws.Range("E2").Font.Color = 15570276

    # Latest Handback File (F2) - mirrors Latest Handoff File (C2) hyperlink
    $ws.Hyperlinks.Add($ws.Range("F2"), $lang.XlfAddress, "", "", $lang.XlfName)
    $ws.Range("F2").Font.Underline = $true
    $ws.Range("F2").Font.Color = 15570276

    # Latest Handback DateTime (G2)
    $ws.Range("G2").Value = $lang.HandbackTime
}
